$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

$sh.Left = 41.51771926879883
$sh.Top = 240.0452880859375
$sh.Width = 442.62994384765625
$sh.Height = 99.49606323242188

$tr = $sh.TextFrame.TextRange
$run = $tr.Runs(7)
$run.Text = "https://github.com/parkzaewoo-ops/reinforcement-PeRL-assignment/"
